# Update the F-column (attendance/visitor count) figures on the
# "展览" (Exhibition) and "全部类型" (All types) sheets.
# Both sheets share identical row layouts for these entries.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value
$updates = @{
    3  = 3395
    6  = 217
    8  = 1658
    16 = 6
    17 = 34
    23 = 63
    24 = 48
    26 = 407
    27 = 268
    28 = 118
    29 = 45
    32 = 454
    33 = 2321
    37 = 575
    39 = 433
    40 = 240
    41 = 358
    43 = 545
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
